$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-22 22:48:52"
$ws.Range("O2").Value = "5.9 °C"
$ws.Range("E3").Value = "2026-02-22 22:48:54"
$ws.Range("H3").Value = "'34%"
$ws.Range("E4").Value = "2026-02-22 22:48:57"
$ws.Range("K4").Value = "14.5 MJ/m2"
$ws.Range("O4").Value = "12.1 °C"
$ws.Range("E5").Value = "2026-02-22 22:48:59"
$ws.Range("N5").Value = "1.8 °C 22:27 TU"
$ws.Range("E6").Value = "2026-02-22 22:49:02"
$ws.Range("J6").Value = "1027.0 hPa"
$ws.Range("E7").Value = "2026-02-22 22:49:05"
$ws.Range("E8").Value = "2026-02-22 22:49:08"
$ws.Range("E9").Value = "2026-02-22 22:49:10"
$ws.Range("H9").Value = "'78%"
$ws.Range("O9").Value = "10.8 °C"
$ws.Range("E10").Value = "2026-02-22 22:49:13"
$ws.Range("H10").Value = "'82%"
$ws.Range("O10").Value = "9.8 °C"
$ws.Range("E11").Value = "2026-02-22 22:49:15"
$ws.Range("H11").Value = "'66%"
$ws.Range("O11").Value = "8.5 °C"
$ws.Range("E12").Value = "2026-02-22 22:49:18"
$ws.Range("H12").Value = "'88%"
$ws.Range("E13").Value = "2026-02-22 22:49:20"
$ws.Range("H13").Value = "'62%"
$ws.Range("E14").Value = "2026-02-22 22:49:23"
$ws.Range("H14").Value = "'75%"
$ws.Range("O14").Value = "11.7 °C"
$ws.Range("E15").Value = "2026-02-22 22:49:26"
$ws.Range("E16").Value = "2026-02-22 22:49:28"
$ws.Range("E17").Value = "2026-02-22 22:49:31"
$ws.Range("N17").Value = "7.6 °C 22:29 TU"
$ws.Range("E18").Value = "2026-02-22 22:49:34"
$ws.Range("H18").Value = "'78%"
$ws.Range("O18").Value = "9.9 °C"
$ws.Range("E19").Value = "2026-02-22 22:49:37"
$ws.Range("E20").Value = "2026-02-22 22:49:39"
$ws.Range("O20").Value = "4.2 °C"
$ws.Range("E21").Value = "2026-02-22 22:49:42"
$ws.Range("E22").Value = "2026-02-22 22:49:45"
$ws.Range("E23").Value = "2026-02-22 22:49:47"
$ws.Range("E24").Value = "2026-02-22 22:49:50"
$ws.Range("O24").Value = "7.6 °C"
$ws.Range("E25").Value = "2026-02-22 22:49:53"
$ws.Range("O25").Value = "6.9 °C"
$ws.Range("E26").Value = "2026-02-22 22:49:55"
$ws.Range("J26").Value = "1026.2 hPa"
$ws.Range("K26").Value = "15.5 MJ/m2"
$ws.Range("O26").Value = "11.0 °C"
$ws.Range("E27").Value = "2026-02-22 22:49:58"
$ws.Range("E28").Value = "2026-02-22 22:50:01"
$ws.Range("E29").Value = "2026-02-22 22:50:04"
$ws.Range("H29").Value = "'84%"
$ws.Range("O29").Value = "9.6 °C"
$ws.Range("E30").Value = "2026-02-22 22:50:06"
$ws.Range("O30").Value = "12.0 °C"
$ws.Range("E31").Value = "2026-02-22 22:50:09"
$ws.Range("L31").Value = "41.4 km/h - 328º 22:29 TU"
$ws.Range("O31").Value = "14.8 °C"
$ws.Range("E32").Value = "2026-02-22 22:50:12"
$ws.Range("O32").Value = "5.6 °C"
$ws.Range("E33").Value = "2026-02-22 22:50:15"
$ws.Range("O33").Value = "8.2 °C"
$ws.Range("E34").Value = "2026-02-22 22:50:17"
$ws.Range("E35").Value = "2026-02-22 22:50:20"
$ws.Range("J35").Value = "1028.5 hPa"
$ws.Range("O35").Value = "11.1 °C"
$ws.Range("E36").Value = "2026-02-22 22:50:23"
$ws.Range("E37").Value = "2026-02-22 22:50:25"
$ws.Range("H37").Value = "'66%"
$ws.Range("J37").Value = "1029.9 hPa"
$ws.Range("O37").Value = "7.8 °C"
$ws.Range("E38").Value = "2026-02-22 22:50:28"
$ws.Range("E39").Value = "2026-02-22 22:50:30"
$ws.Range("E40").Value = "2026-02-22 22:50:33"
$ws.Range("H40").Value = "'58%"
$ws.Range("J40").Value = "1029.2 hPa"
$ws.Range("O40").Value = "9.6 °C"
$ws.Range("E41").Value = "2026-02-22 22:50:35"
$ws.Range("H41").Value = "'79%"
$ws.Range("O41").Value = "11.1 °C"
$ws.Range("E42").Value = "2026-02-22 22:50:38"
$ws.Range("E43").Value = "2026-02-22 22:50:41"
$ws.Range("O43").Value = "9.3 °C"
$ws.Range("E44").Value = "2026-02-22 22:50:43"
$ws.Range("K44").Value = "15.6 MJ/m2"
$ws.Range("L44").Value = "32.4 km/h - 42º 22:21 TU"
$ws.Range("E45").Value = "2026-02-22 22:50:46"
$ws.Range("O45").Value = "8.5 °C"
$ws.Range("E46").Value = "2026-02-22 22:50:48"
$ws.Range("J46").Value = "1029.4 hPa"
